# Update InsideBet Data: Automatizado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Milan ---
$ws.Range("C3").Value = 25
$ws.Range("E3").Value = 9
$ws.Range("G3").Value = 41
$ws.Range("H3").Value = 19
$ws.Range("J3").Value = 54
$ws.Range("K3").Value = 2.16
$ws.Range("L3").Value = "W D W W D"
$ws.Range("M3").Value = 72830
$ws.Range("N3").Value = "Christian Pulisic, Rafael Leão - 8"

# --- Row 7: was Atalanta, now Como (team order swap in the table) ---
$ws.Range("B7").Value = "Como"
$ws.Range("G7").Value = 39
$ws.Range("H7").Value = 19
$ws.Range("I7").Value = 20
$ws.Range("L7").Value = "W W D L D"
$ws.Range("M7").Value = 11409
$ws.Range("N7").Value = "Nicolás Paz - 9"
$ws.Range("O7").Value = "Jean Butez"

# --- Row 8: was Como, now Atalanta ---
$ws.Range("B8").Value = "Atalanta"
$ws.Range("C8").Value = 25
$ws.Range("E8").Value = 9
$ws.Range("G8").Value = 34
$ws.Range("H8").Value = 21
$ws.Range("I8").Value = 13
$ws.Range("J8").Value = 42
$ws.Range("K8").Value = 1.68
$ws.Range("L8").Value = "D W D W W"
$ws.Range("M8").Value = 22163
$ws.Range("N8").Value = "Nikola Krstović - 7"
$ws.Range("O8").Value = "Marco Carnesecchi"

# --- Row 14: Cagliari attendance update ---
$ws.Range("M14").Value = 16023

# --- Row 15: Torino attendance update ---
$ws.Range("M15").Value = 19108
